# Add data for 2024-12-24 -- refresh 2024 (and a couple 2023 corrections) totals
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7438
$ws.Range("K3").Value = 7700
$ws.Range("J4").Value = 1848
$ws.Range("K4").Value = 1619
$ws.Range("K6").Value = 8557
$ws.Range("J7").Value = 29317
$ws.Range("K7").Value = 25861

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 224
$ws.Range("K7").Value = 780
$ws.Range("K8").Value = 1680
$ws.Range("K9").Value = 125
$ws.Range("K11").Value = 464
$ws.Range("K19").Value = 750
$ws.Range("K23").Value = 261
$ws.Range("K29").Value = 1425
$ws.Range("K33").Value = 1093
$ws.Range("K36").Value = 328
$ws.Range("K37").Value = 864
$ws.Range("K41").Value = 175
$ws.Range("K42").Value = 951
$ws.Range("K44").Value = 210
$ws.Range("K45").Value = 36
$ws.Range("K47").Value = 177
$ws.Range("K51").Value = 326
$ws.Range("K55").Value = 282
$ws.Range("K56").Value = 30
$ws.Range("J63").Value = 128
$ws.Range("K63").Value = 79
$ws.Range("K65").Value = 608
$ws.Range("K67").Value = 1011
$ws.Range("K76").Value = 357
$ws.Range("K78").Value = 318
$ws.Range("K79").Value = 634
$ws.Range("K83").Value = 545
$ws.Range("K84").Value = 210
$ws.Range("K85").Value = 1185
$ws.Range("K88").Value = 276
$ws.Range("K89").Value = 387
$ws.Range("K94").Value = 344
$ws.Range("K99").Value = 435
$ws.Range("J101").Value = 29317
$ws.Range("K101").Value = 25861

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K6").Value = 218
$ws.Range("K7").Value = 780

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K6").Value = 161
$ws.Range("K7").Value = 464

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 116
$ws.Range("K7").Value = 387

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 390
$ws.Range("K3").Value = 410
$ws.Range("K7").Value = 1185

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 181
$ws.Range("K4").Value = 38

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K6").Value = 564
$ws.Range("K7").Value = 1680

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 189
$ws.Range("K6").Value = 128
$ws.Range("K7").Value = 545

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K6").Value = 349
$ws.Range("K7").Value = 1093

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 286
$ws.Range("K7").Value = 864

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K6").Value = 225
$ws.Range("K7").Value = 608

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 435

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 285
$ws.Range("K7").Value = 1011

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 210

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 403
$ws.Range("K3").Value = 506
$ws.Range("K4").Value = 66
$ws.Range("K7").Value = 1425

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 223
$ws.Range("K7").Value = 750

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 57
$ws.Range("K7").Value = 210

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K4").Value = 29
$ws.Range("K6").Value = 179
$ws.Range("K7").Value = 357

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 175

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 259
$ws.Range("K6").Value = 352
$ws.Range("K7").Value = 951

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K6").Value = 107
$ws.Range("K7").Value = 318

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 82
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 282

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 72
$ws.Range("K7").Value = 261

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K4").Value = 42
$ws.Range("K7").Value = 634

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 122
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 328

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 82
$ws.Range("K4").Value = 28
$ws.Range("K7").Value = 344

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 177

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 224

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 87
$ws.Range("K7").Value = 276

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 88
$ws.Range("K7").Value = 326

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 30
